$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the typo'd e-mail address: escig15@gmail.com -> efeig15@gmail.com
# ------------------------------------------------------------------
$d.Content.Find.Execute("escig15@gmail.com", $false, $false, $false, $false, $false,
                         $true, 1, $false, "efeig15@gmail.com", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Fill in the two blanks ("___________") in the paragraph that asks
#    who to contact with questions, splitting each run of 11
#    underscores into: "_" + <filled-in value, underlined> + remaining
#    underscores, keeping the yellow highlight throughout.
# ------------------------------------------------------------------
$paraStart = -1
$paraEnd = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*___________*" -and $t -like "*בטלפון*") {
        $paraStart = $p.Range.Start
        $paraEnd = $p.Range.End
        break
    }
}

if ($paraStart -ge 0) {

    # ---- first blank: contact person's name ----
    $searchRange = $d.Range($paraStart, $paraEnd)
    $searchRange.Find.Execute("___________", $false, $false, $false, $false, $false,
                               $true, 1, $false, "", 0) | Out-Null

    # Shrink the run down to a single leading underscore.
    $searchRange.Text = "_"

    # Insert the name right after the leading underscore.
    $nameRange = $d.Range($searchRange.End, $searchRange.End)
    $nameRange.InsertAfter("אלישר פייג")

    # Insert the trailing underscores right after the name.
    $tailRange = $d.Range($nameRange.End, $nameRange.End)
    $tailRange.InsertAfter("__")

    # Underline only the name - done last so neighboring runs don't inherit it.
    $nameRange.Font.Underline = 1

    # ---- second blank: contact person's phone number ----
    $searchRange2 = $d.Range($tailRange.End, $paraEnd + 20)
    $searchRange2.Find.Execute("___________", $false, $false, $false, $false, $false,
                                $true, 1, $false, "", 0) | Out-Null

    $searchRange2.Text = "_"

    $phoneRange = $d.Range($searchRange2.End, $searchRange2.End)
    $phoneRange.InsertAfter("058-7272372")

    $tailRange2 = $d.Range($phoneRange.End, $phoneRange.End)
    $tailRange2.InsertAfter("_")

    $phoneRange.Font.Underline = 1
}
